$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.480.60"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "2.164.78"
$ws.Range("E3").Value = "  +3.00%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'227.88"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("D7").Value = "'64.03"
$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.396"
$ws.Range("E9").Value = "  +2.26%  "

$ws.Range("D10").Value = "'0.0856"
$ws.Range("E10").Value = "  +1.71%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").Value = "'16.19"
$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("D13").Value = "2.485.83"
$ws.Range("E13").Value = "  +3.13%  "

$ws.Range("D14").Value = "'22.14"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "'0.814"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "'5.53"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").Value = "2.166.30"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("D18").Value = "39.501.49"
$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("D19").Value = "'71.98"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").Value = "0.0₃0851"
$ws.Range("E21").Value = "  +1.35%  "

$ws.Range("D22").Value = "'229.61"
$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -0.76%  "

$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'172.21"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.61"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("D29").Value = "'1.44"
$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("D30").Value = "'19.91"
$ws.Range("E30").Value = "  +2.91%  "

$ws.Range("D31").Value = "'2.60"
$ws.Range("E31").Value = "  +5.01%  "

$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "'4.63"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("E35").Value = "  -0.77%  "

$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").Value = "'2.44"
$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("D38").Value = "'3.57"
$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "'103.41"
$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "'17.82"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").Value = "1.525.53"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("E44").Value = "  +3.76%  "

$ws.Range("D45").Value = "'0.0932"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  +5.68%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'4.26"
$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("D49").Value = "'7.76"
$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("D50").Value = "2.369.42"
$ws.Range("E50").Value = "  +3.37%  "

$ws.Range("E51").Value = "  -0.43%  "
